$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    # Force the cell to stay a text cell so numeric-looking strings
    # (e.g. "536.70", "1.00") keep their exact formatting instead of
    # being coerced to a Double and losing trailing zeros / precision.
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '59.748.68'
$ws.Range('E2').Value = '  +3.31%  '
$ws.Range('D3').Value = '3.189.17'
$ws.Range('E3').Value = '  +2.31%  '
$ws.Range('E4').Value = '  +0.05%  '
Set-TextValue $ws 'D5' '536.70'
$ws.Range('E5').Value = '  +0.60%  '
Set-TextValue $ws 'D6' '144.97'
$ws.Range('E6').Value = '  +4.66%  '
$ws.Range('E7').Value = '  -0.09%  '
Set-TextValue $ws 'D8' '0.520'
$ws.Range('E8').Value = '  +4.90%  '
$ws.Range('E9').Value = '  -0.96%  '
$ws.Range('E10').Value = '  +5.13%  '
Set-TextValue $ws 'D11' '0.430'
$ws.Range('E11').Value = '  +4.00%  '
$ws.Range('D12').Value = '3.738.53'
$ws.Range('E12').Value = '  +2.30%  '
$ws.Range('E13').Value = '  +0.07%  '
Set-TextValue $ws 'D14' '25.98'
$ws.Range('E14').Value = '  +0.85%  '
$ws.Range('E15').Value = '  +4.24%  '
$ws.Range('D16').Value = '59.824.55'
$ws.Range('E16').Value = '  +3.32%  '
$ws.Range('D17').Value = '3.215.06'
$ws.Range('E17').Value = '  +3.47%  '
$ws.Range('E18').Value = '  +0.85%  '
Set-TextValue $ws 'D19' '13.03'
$ws.Range('E19').Value = '  +1.44%  '
Set-TextValue $ws 'D20' '8.23'
$ws.Range('E20').Value = '  +1.40%  '
Set-TextValue $ws 'D21' '380.52'
$ws.Range('E21').Value = '  +1.92%  '
$ws.Range('E22').Value = '  -0.10%  '
Set-TextValue $ws 'D23' '0.530'
$ws.Range('E23').Value = '  +4.29%  '
Set-TextValue $ws 'D24' '70.16'
$ws.Range('E24').Value = '  +1.19%  '
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws 'D25' '0.172'
$ws.Range('E25').Value = '  +3.01%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws 'D26' '8.85'
$ws.Range('E26').Value = '  +16.77%  '
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('D28').Value = '0.0₃0906'
$ws.Range('E28').Value = '  +2.19%  '
$ws.Range('E29').Value = '  +2.34%  '
Set-TextValue $ws 'D30' '22.42'
$ws.Range('E30').Value = '  +4.16%  '
Set-TextValue $ws 'D31' '6.17'
$ws.Range('E31').Value = '  +0.18%  '
Set-TextValue $ws 'D32' '5.40'
$ws.Range('E32').Value = '  +4.67%  '
Set-TextValue $ws 'D33' '1.20'
$ws.Range('E33').Value = '  +1.83%  '
Set-TextValue $ws 'D34' '6.45'
$ws.Range('E34').Value = '  +4.52%  '
Set-TextValue $ws 'D35' '157.01'
$ws.Range('E35').Value = '  -2.12%  '
$ws.Range('E36').Value = '  +3.77%  '
$ws.Range('B37').Value = 'EnergySwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D37' '25.68'
$ws.Range('E37').Value = '  +0.48%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '2.752.21'
$ws.Range('E38').Value = '  +7.68%  '
Set-TextValue $ws 'D39' '0.0711'
$ws.Range('E39').Value = '  +6.00%  '
Set-TextValue $ws 'D41' '4.29'
$ws.Range('E41').Value = '  +3.32%  '
Set-TextValue $ws 'D42' '0.726'
$ws.Range('E42').Value = '  +4.30%  '
Set-TextValue $ws 'D43' '39.53'
$ws.Range('E43').Value = '  +3.01%  '
$ws.Range('E44').Value = '  +6.89%  '
$ws.Range('D45').Value = '3.234.41'
$ws.Range('E45').Value = '  +2.44%  '
Set-TextValue $ws 'D46' '1.00'
$ws.Range('E46').Value = '  +2.07%  '
Set-TextValue $ws 'D47' '6.19'
$ws.Range('E47').Value = '  +0.75%  '
Set-TextValue $ws 'D48' '0.100'
$ws.Range('E48').Value = '  +6.15%  '
Set-TextValue $ws 'D49' '20.54'
$ws.Range('E49').Value = '  +3.02%  '
Set-TextValue $ws 'D50' '0.775'
$ws.Range('E50').Value = '  +3.74%  '
$ws.Range('E51').Value = '  +0.04%  '

Write-Output "Applied all changes"
